{"js": "// 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", \"Replace\");\n}\n\n// 2. Split the mailing-address line (\"929 Story Road, San Jose CA 95122\")\n//    that sits right under the addressee's name into two separate\n//    paragraphs: \"929 Story Road\" and \"San Jose, CA 95122\".\n//    (There is an identical string inside the PROPERTY ADDRESS table\n//    cell further down which must stay untouched, so we match on the\n//    body paragraph collection and confirm it is NOT inside a table.)\nconst paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nconst candidateIndexes = [];\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"929 Story Road, San Jose CA 95122\") {\n    paras.items[i].load(\"parentTableOrNullObject\");\n    candidateIndexes.push(i);\n  }\n}\nawait context.sync();\n\nlet addressIndex = -1;\nfor (const i of candidateIndexes) {\n  if (paras.items[i].parentTableOrNullObject.isNullObject) {\n    addressIndex = i;\n    break;\n  }\n}\n\nif (addressIndex !== -1) {\n  const addressPara = paras.items[addressIndex];\n  addressPara.insertText(\"929 Story Road\", \"Replace\");\n  addressPara.insertParagraph(\"San Jose, CA 95122\", \"After\");\n  await context.sync();\n}\n\n// 3. Remove the now-redundant blank \"No Spacing\" paragraph that used to\n//    follow the \"...Board of Directors\" line.\nconst paras2 = context.document.body.paragraphs;\nparas2.load(\"text,style\");\nawait context.sync();\n\nlet blankIndex = -1;\nfor (let i = 1; i < paras2.items.length; i++) {\n  if (\n    paras2.items[i].text === \"\" &&\n    paras2.items[i].style === \"No Spacing\" &&\n    paras2.items[i - 1].text.indexOf(\"Board of Directors\") !== -1\n  ) {\n    blankIndex = i;\n    break;\n  }\n}\n\nif (blankIndex !== -1) {\n  paras2.items[blankIndex].delete();\n  await context.sync();\n}\n", "ps1": "# 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$d = $word.ActiveDocument\n\n$count0 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count0; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = ($p.Range.Text -replace \"[\\r\\a\\x07]+$\", \"\")\n    if ($txt -eq \"September 19, 2025\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2. Split the mailing-address line (\"929 Story Road, San Jose CA 95122\")\n#    that sits right under the addressee's name into two separate\n#    paragraphs: \"929 Story Road\" and \"San Jose, CA 95122\".\n#    (There is an identical string inside the PROPERTY ADDRESS table\n#    cell further down which must stay untouched, so confirm the\n#    paragraph is NOT inside a table before touching it.)\n$count = $d.Paragraphs.Count\n$addressIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = ($p.Range.Text -replace \"[\\r\\a\\x07]+$\", \"\")\n    if ($txt -eq \"929 Story Road, San Jose CA 95122\") {\n        $inTable = $p.Range.Information(12)\n        if (-not $inTable) {\n            $addressIndex = $i\n            break\n        }\n    }\n}\n\nif ($addressIndex -gt 0) {\n    $addrPara = $d.Paragraphs.Item($addressIndex)\n    $addrRange = $addrPara.Range\n    $addrRange.Text = \"929 Story Road\"\n    $addrRange.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($addressIndex + 1)\n    $newPara.Range.Text = \"San Jose, CA 95122\"\n}\n\n# 3. Remove the now-redundant blank \"No Spacing\" paragraph that used to\n#    follow the \"...Board of Directors\" line.\n$count2 = $d.Paragraphs.Count\n$blankIndex = -1\nfor ($i = 1; $i -le $count2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = ($p.Range.Text -replace \"[\\r\\a\\x07]+$\", \"\")\n    if ($txt -like \"*Board of Directors*\") {\n        $next = $d.Paragraphs.Item($i + 1)\n        $nextTxt = ($next.Range.Text -replace \"[\\r\\a\\x07]+$\", \"\")\n        if ($nextTxt -eq \"\" -and $next.Style.NameLocal -eq \"No Spacing\") {\n            $blankIndex = $i + 1\n        }\n        break\n    }\n}\n\nif ($blankIndex -gt 0) {\n    $d.Paragraphs.Item($blankIndex).Range.Delete()\n}\n"}
